$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Make "Repayment schedule" the active sheet/tab (was "Transactions" before).
$ws.Activate()

# Insert a new blank column before column N (14) - "Late"/"heading"/"Outstanding"
# columns shift one place to the right (N->O, O->P, P->Q).
$ws.Columns.Item(14).Insert()

# The freshly inserted column picks up the width of the column to its left.
$ws.Columns.Item(14).ColumnWidth = 9.8

# Update the selection on the Repayment schedule sheet.
$ws.Range("J19").Select() | Out-Null
